$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of row 5 (A5:V5) while preserving cell formatting/styles.
$ws.Range("A5:V5").ClearContents()
